$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 140 ("「字は手による声なり」" entry) entirely; this shifts all
# subsequent rows up by one, matching the author's removal of that post.
$ws.Rows.Item(140).Delete()
